$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells C1/D1 - bold + centered, matching existing header style
$ws.Range("C1").Value = "ci_lower"
$ws.Range("D1").Value = "ci_upper"
$ws.Range("C1:D1").Font.Bold = $true
$ws.Range("C1:D1").HorizontalAlignment = -4108

# Confidence interval data, rows 2-74
$cLower = @(
4366.965832188154, 4508.126607641642, 4645.433736053693, 4777.720911464151, 4906.53770436577, 5032.159999355847, 5154.049344640302, 5265.567592250796, 5368.425421184447, 5467.176787954823, 5558.629717911714, 5645.276323311718, 5723.906436550624, 5795.065431066633, 5858.602183710405, 5917.240606866138, 5969.307698691387, 6014.52576631375, 6052.516590140955, 6083.288829081212, 6107.813662424499, 6123.600775073039, 6129.775413652443, 6132.825830920911, 6127.515532219901, 6113.162105155106, 6094.781858845976, 6065.834615139046, 6033.026694646488, 5990.435382203454, 5940.202696437649, 5880.980029535578, 5817.508370612177, 5740.570670306622, 5659.418728124768, 5572.437573548558, 5477.311239722902, 5376.855004098755, 5268.507859760133, 5155.477457837166, 5035.753459447461, 4911.231281397589, 4783.879106697284, 4653.978962883332, 4517.398705864073, 4379.821443723634, 4237.529370973304, 4095.121500884416, 3951.927881194395, 3807.695287209804, 3660.262121882606, 3518.40935375347, 3377.783866406036, 3235.873570247929, 3091.673001976354, 2953.387245452378, 2820.758149111723, 2688.560932784423, 2557.283494466538, 2428.082791700785, 2303.207102730948, 2181.621909116167, 2063.483324109895, 1949.034654529953, 1838.370744479844, 1731.542083613694, 1628.618691007788, 1530.452224891158, 1436.337484100249, 1345.401498578323, 1257.866827159485, 1173.943653529322, 1094.450888304071
)
$dUpper = @(
4646.843522565296, 4767.578463546169, 4884.180434324798, 4997.992756061139, 5107.296298882908, 5219.137573467197, 5326.296545645789, 5428.373404314416, 5530.321634467779, 5629.797155178417, 5727.437722801433, 5821.232617327908, 5905.392248315897, 5986.599023749786, 6060.060538332566, 6128.972766850394, 6190.064305820034, 6240.578959322797, 6285.645686796362, 6323.561494373513, 6355.184824767187, 6376.143366005647, 6386.949078062442, 6388.75950856598, 6384.336904020161, 6370.956918434043, 6350.48151974775, 6318.325846979027, 6280.597711659953, 6236.781855929858, 6186.44482446449, 6129.952474861789, 6067.176661554247, 6001.47661580976, 5928.244190229791, 5841.340037095031, 5757.277105187518, 5671.156893245293, 5578.489959259067, 5474.246886289774, 5372.138126142343, 5266.354811254862, 5161.117767222458, 5055.556910613967, 4940.395465425464, 4823.855237864888, 4704.981906161242, 4587.149893900279, 4468.599779314331, 4348.455774514429, 4228.66856901791, 4107.533888167538, 3985.21101142678, 3860.176063486586, 3734.889724627717, 3609.717122585584, 3485.00739732452, 3362.738226865609, 3242.18005579247, 3122.615396834072, 3004.257552166029, 2885.843424173981, 2767.261314391807, 2651.576947966274, 2540.891705808763, 2429.240129471683, 2319.285399470349, 2211.789123072723, 2106.887323205922, 2005.711331884272, 1907.470936262407, 1809.771369925798, 1715.798323936552
)

for ($i = 0; $i -lt $cLower.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $cLower[$i]
    $ws.Cells.Item($row, 4).Value = $dUpper[$i]
}

Write-Output "done"